$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "List1"

# New data values (plain numbers, no formulas, no special style)
$data = @(
    @(1, 3),
    @(3, 4),
    @(7, 2),
    @(13, 2),
    @(17, 2),
    @(22, 2),
    @(28, 2),
    @(32, 2),
    @(43, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove old row 10 contents (since new range is only A1:B9)
$ws.Range("A10:B10").Clear()

# Clear any styling previously applied to column B (s="1")
$ws.Range("A1:B10").ClearFormats()

# Update selection to A2
$ws.Range("A2").Select()
